$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values (row 3 and row 4) for the columns that need to be swapped
$a3 = $ws.Range("A3").Value2
$a4 = $ws.Range("A4").Value2
$i3 = $ws.Range("I3").Value2
$i4 = $ws.Range("I4").Value2
$q3 = $ws.Range("Q3").Value2
$q4 = $ws.Range("Q4").Value2
$r3 = $ws.Range("R3").Value2
$r4 = $ws.Range("R4").Value2

# Column A (Id) - numeric swap
$ws.Range("A3").Value2 = $a4
$ws.Range("A4").Value2 = $a3

# Column I (Antal) - keep as text, force text format so numeric-looking strings
# are not auto-converted to numbers by Excel
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I3").Value2 = [string]$i4
$ws.Range("I4").Value2 = [string]$i3

# Column Q (Ost) - numeric swap
$ws.Range("Q3").Value2 = $q4
$ws.Range("Q4").Value2 = $q3

# Column R (Nord) - numeric swap
$ws.Range("R3").Value2 = $r4
$ws.Range("R4").Value2 = $r3
